$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.512.31"
$ws.Range("E2").Value = "  +0.54%  "

$ws.Range("D3").Value = "2.240.59"
$ws.Range("E3").Value = "  -0.42%  "

$ws.Range("E4").Value = "  +0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.97%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("E8").Value = "  +0.18%  "

$ws.Range("E9").Value = "  -2.32%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0800"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.28%  "

$ws.Range("E12").Value = "  -1.92%  "

$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").Value = "2.582.44"
$ws.Range("E14").Value = "  -0.50%  "

$ws.Range("D15").Value = "2.331.16"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.830"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.96%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.94%  "

$ws.Range("D18").Value = "44.313.27"
$ws.Range("E18").Value = "  +0.42%  "

$ws.Range("D19").Value = "0.0₃0937"
$ws.Range("E19").Value = "  -3.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.17"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.70%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.32%  "

$ws.Range("E26").Value = "  -0.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.10%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.49%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "37.03"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.30%  "

$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.65%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.85%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.49%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0783"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.53%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.14"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.79%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.11%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("E38").Value = "  -1.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.75%  "

$ws.Range("E42").Value = "  -1.78%  "

$ws.Range("E43").Value = "  +0.10%  "

$ws.Range("D44").Value = "1.804.19"
$ws.Range("E44").Value = "  +2.50%  "

$ws.Range("E45").Value = "  +12.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "80.69"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "

$ws.Range("E47").Value = "  -3.18%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "97.94"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "

$ws.Range("E49").Value = "  -3.08%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.29"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.37%  "

$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.15%  "
